# feat: add 2022-Q3 data
#
# - "总计" sheet: the Q2 summary row becomes the new Q3 summary row (same
#   market value that used to live in the single data row moves down into a
#   freshly appended row for 2022-Q2, and the original row is overwritten
#   with the new 2022-Q3 numbers).
# - The existing "2022-Q2" detail sheet is duplicated so the untouched
#   Q2 fund-holding detail survives under its original name, while the
#   original sheet (same sheetId/rId) is renamed to "2022-Q3" and repopulated
#   with the new quarter's fund-holding detail.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet: shift the existing Q2 row down to row 3, and put the
#    new Q3 figures on row 2.
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item(1)

# Give the new row 3 the same style as row 2's A-cell (bold/boxed), then
# write the old Q2 values into it.
$totals.Range("A2").Copy() | Out-Null
$totals.Range("A3").PasteSpecial(-4122) | Out-Null

$totals.Range("A3").Value = 1
$totals.Range("B3").Value = "2022-Q2"
$totals.Range("C3").Value = 2
$totals.Range("D3").Value = 0.19

# Overwrite row 2 with the new Q3 figures.
$totals.Range("B2").Value = "2022-Q3"
$totals.Range("D2").Value = 0.12

# ---------------------------------------------------------------------
# 2) Duplicate the current "2022-Q2" sheet so the original detail data is
#    preserved unchanged under its original name, positioned right after
#    the source sheet.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item(2)
$q2.Copy($null, $q2) | Out-Null

# Free up the "2022-Q2" name on the original sheet before renaming the
# copy to it.
$q2.Name = "2022-Q3"
$wb.Worksheets.Item(3).Name = "2022-Q2"

# ---------------------------------------------------------------------
# 3) Re-style the (now renamed) "2022-Q3" sheet's header row and A-column
#    to match the bold/boxed look used on the "总计" sheet, then replace
#    the fund metrics with the new quarter's numbers. Fund code/name stay
#    the same.
# ---------------------------------------------------------------------
$totals.Range("B1:D1").Copy() | Out-Null
$q2.Range("B1:H1").PasteSpecial(-4122) | Out-Null

$totals.Range("A2").Copy() | Out-Null
$q2.Range("A2:A3").PasteSpecial(-4122) | Out-Null

# Row 2 (014275 / 易方达北交所精选两年定开混合A)
$q2.Range("D2").Value = "'3.58"
$q2.Range("E2").Value = "'61.75"
$q2.Range("F2").Value = "'2.65"
$q2.Range("G2").Value = "'0.0949"
$q2.Range("H2").Value = 7

# Row 3 (014276 / 易方达北交所精选两年定开混合C)
$q2.Range("D3").Value = "'0.92"
$q2.Range("E3").Value = "'61.75"
$q2.Range("F3").Value = "'2.65"
$q2.Range("G3").Value = "'0.0244"
$q2.Range("H3").Value = 7
